$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$akCol = 37  # Column AK
$lastRow = $ws.Cells.Item($ws.Rows.Count, $akCol).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $akCol)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val / 10000
    }
}
